$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add a new daily-report block (16/10/2014) by cloning the previous
# day's block (rows 49-57, 14/10/2014) into rows 58-66, then updating
# the date and the free-text content cells.
# ------------------------------------------------------------------
$src = $ws.Range("B49:C57")
$dst = $ws.Range("B58")
$src.Copy($dst)

# New date: 16/10/2014 (serial 41928)
$ws.Range("B58").Value = 41928

# "Ke hoach" (plan) text for the new day
$ws.Range("C59").Value = "``- Công việc 1:Buid giao diện màn hình Tìm Kiếm và kết nối API cho màn hình này.`n``- Công việc 2: Buid giao diện màn hình Công việc theo dõi`n``- Công việc 2: Buid giao diện màn hình Công việc hoàn thành"

# "Ket qua dat duoc" (result) text - needs quote-prefix formatting preserved
$ws.Range("C60").Value = "'- Công việc 1: 90%`n'- Công việc 2: hoàn thành`n'- Công việc 2: hoàn thành"

# C61 (Trang thai / status) and C63 (Giai quyet van de header) already
# carry over the correct, reused text from the copy above - no change
# needed there.

# "Van de gap phai" (issue) text - needs quote-prefix formatting preserved
$ws.Range("C62").Value = "'- Vấn đề 1: Chưa có API Tìm kiếm"

# "Giai quyet van de" (solution) text - needs quote-prefix formatting preserved
$ws.Range("C64").Value = "'- Vấn đề 1: Yêu cầu bộ phận thiết kế API hoàn thành ."

# Second "Giai quyet van de:" label, this time re-typed with a leading
# quote (re-using the existing shared string) which keeps its
# quote-prefix cell format.
$ws.Range("C65").Value = "'Giải quyết vấn đề:"

# "Ke hoach ngay mai" (tomorrow's plan) text
$ws.Range("C66").Value = "Hoàn thành màn hình Tìm kiếm`nvà kết nối API cho các màn hình còn lại"

# Row heights for the wrapped, multi-line cells
$ws.Rows.Item(59).RowHeight = 45
$ws.Rows.Item(60).RowHeight = 45
$ws.Rows.Item(66).RowHeight = 30

# Update the view: scroll down and select B58, same as the source edit
$ws.Range("A56").Select()
$excel.ActiveWindow.ScrollRow = 56
$ws.Range("B58").Select()
